$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells stay as text, since some values
# look numeric and Excel would otherwise coerce them to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.036.82'
$ws.Range("D3").Value = '1.645.89'
$ws.Range("D5").Value = '216.77'
$ws.Range("D9").Value = '0.0642'
$ws.Range("D10").Value = '19.69'
$ws.Range("D11").Value = '0.0795'
$ws.Range("D13").Value = '1.871.80'
$ws.Range("D14").Value = '1.655.71'
$ws.Range("D17").Value = '63.16'
$ws.Range("D18").Value = '26.007.48'
$ws.Range("D20").Value = '193.50'
$ws.Range("D22").Value = '9.96'
$ws.Range("D24").Value = '0.134'
$ws.Range("D26").Value = '144.61'
$ws.Range("D28").Value = '6.93'
$ws.Range("D29").Value = '15.56'
$ws.Range("D32").Value = '3.30'
$ws.Range("D37").Value = '1.134.35'
$ws.Range("D38").Value = '0.543'
$ws.Range("D42").Value = '99.57'
$ws.Range("D43").Value = '0.798'
$ws.Range("D44").Value = '1.780.72'
$ws.Range("D46").Value = '56.84'

# Restore the default (Normal) style on the price column so no
# stray number formatting is left behind on the cells.
$priceRange.Style = "Normal"

# Update the Volume(1h) column (E) percentages.
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").Value = '  +0.64%  '
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("E6").Value = '  +1.60%  '
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("E8").Value = '  +0.98%  '
$ws.Range("E9").Value = '  +1.81%  '
$ws.Range("E10").Value = '  +0.55%  '
$ws.Range("E11").Value = '  +1.07%  '
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("E13").Value = '  +1.04%  '
$ws.Range("E14").Value = '  +2.15%  '
$ws.Range("E15").Value = '  +0.68%  '
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("E18").Value = '  +0.60%  '
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("E24").Value = '  +8.14%  '
$ws.Range("E25").Value = '  +2.10%  '
$ws.Range("E26").Value = '  +1.42%  '
$ws.Range("E27").Value = '  +0.69%  '
$ws.Range("E28").Value = '  +1.21%  '
$ws.Range("E29").Value = '  +1.07%  '
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("E32").Value = '  -0.38%  '
$ws.Range("E33").Value = '  +1.54%  '
$ws.Range("E34").Value = '  -2.50%  '
$ws.Range("E35").Value = '  +3.00%  '
$ws.Range("E36").Value = '  +0.79%  '
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("E40").Value = '  +1.10%  '
$ws.Range("E41").Value = '  +0.78%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("E45").Value = '  +3.60%  '
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("E48").Value = '  +0.60%  '
$ws.Range("E49").Value = '  +1.66%  '
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("E51").Value = '  +0.40%  '
